# Add a default header to the only section, containing the centered,
# Arial 12pt text "Questionnaire 49" (styled with the built-in "Header"
# paragraph style), so the questionnaire number survives printing.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)   # wdHeaderFooterPrimary

# InsertAfter (rather than setting .Range.Text) mints only the "default"
# header part/reference instead of also minting even-page / first-page
# header+footer parts.
$header.Range.InsertAfter("Questionnaire 49")

# Paragraph-level formatting: built-in Header style, centered.
$header.Range.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Run-level formatting: Arial, 12pt (sz is in half-points -> 24).
# Apply to a range that excludes the trailing paragraph mark so the
# formatting lands on the run only, not on <w:pPr><w:rPr>.
$textRange = $header.Range.Duplicate
$textRange.MoveEnd(1, -1) | Out-Null
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
